# Update "想去人数" (F) and "最低票价" (G) figures on the two sheets that
# carry this event table: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> hashtable of column letter -> new value
$updates = @{
    2  = @{ F = 306; G = 250 }
    3  = @{ F = 11329 }
    4  = @{ F = 10647 }
    5  = @{ F = 600 }
    7  = @{ F = 760 }
    8  = @{ F = 109 }
    9  = @{ F = 32 }
    12 = @{ F = 10525 }
    13 = @{ F = 3251 }
    19 = @{ F = 407 }
    20 = @{ F = 11081 }
    21 = @{ F = 10829 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $address = "$col$row"
            $ws.Range($address).Value = $cols[$col]
        }
    }
}
